# Time Trial table (supervisor view): the "TimeTrial" / "TimeTrialsOptionTime"
# blocks on Sheet1 used to start two columns further right than they needed to
# (columns I:O), leaving columns G:H completely empty. Pull the tables back
# over so they start at column G (G:M) by deleting the two empty spacer
# columns - Excel shifts all the data/formulas/formatting in columns I:O left
# by two and re-points the formulas automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1:H1").EntireColumn.Delete()
